$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers
$ws.Range("H1").Value = "Tempo Heuristica"
$ws.Range("I1").Value = "Tempo Total"

# Row 2 updates
$ws.Range("G2").Value = 0.05084657669067383
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.05084657669067383

# Row 3 updates
$ws.Range("G3").Value = 0.05266118049621582
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0.05266118049621582

# Row 4 updates
$ws.Range("G4").Value = 0.05157780647277832
$ws.Range("H4").Value = 0
$ws.Range("I4").Value = 0.05157780647277832
